$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.609230666666667
$ws.Range("H2").Value = 4.827692
$ws.Range("I2").Value = 0.5482851650894511
$ws.Range("J2").Value = 0.5482851650894512
$ws.Range("M2").Value = 61.156892
$ws.Range("N2").Value = 183.470676
$ws.Range("O2").Value = 0.9308124486389074
$ws.Range("P2").Value = 0.9308124486389074
$ws.Range("Q2").Value = 98.41554608442134
$ws.Range("R2").Value = 885.739914759792
$ws.Range("S2").Value = 0.5103506570692996
$ws.Range("T2").Value = 0.5103506570692997
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.609230666666667
$ws.Range("H3").Value = 4.827692
$ws.Range("I3").Value = 0.5482851650894511
$ws.Range("J3").Value = 0.5482851650894512
$ws.Range("O3").Value = 0.02171808228502914
$ws.Range("P3").Value = 0.02171808228502914
$ws.Range("Q3").Value = 2.296270243391111
$ws.Range("R3").Value = 20.66643219052
$ws.Range("S3").Value = 0.01190770233107349
$ws.Range("T3").Value = 0.01190770233107349
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.609230666666667
$ws.Range("H4").Value = 4.827692
$ws.Range("I4").Value = 0.5482851650894511
$ws.Range("J4").Value = 0.5482851650894512
$ws.Range("M4").Value = 2.00294
$ws.Range("N4").Value = 6.00882
$ws.Range("O4").Value = 0.03048489589491914
$ws.Range("P4").Value = 0.03048489589491914
$ws.Range("Q4").Value = 3.223192471493334
$ws.Range("R4").Value = 29.00873224344
$ws.Range("S4").Value = 0.01671441617848047
$ws.Range("T4").Value = 0.01671441617848048
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.609230666666667
$ws.Range("H5").Value = 4.827692
$ws.Range("I5").Value = 0.5482851650894511
$ws.Range("J5").Value = 0.5482851650894512
$ws.Range("M5").Value = 1.115932333333334
$ws.Range("N5").Value = 3.347797
$ws.Range("O5").Value = 0.01698457318114416
$ws.Range("P5").Value = 0.01698457318114415
$ws.Range("Q5").Value = 1.795792532724889
$ws.Range("R5").Value = 16.162132794524
$ws.Range("S5").Value = 0.009312389510597489
$ws.Range("T5").Value = 0.009312389510597489
# Row 6
$ws.Range("G6").Value = 0.9591736666666666
$ws.Range("I6").Value = 0.3268025542087943
$ws.Range("J6").Value = 0.3268025542087943
$ws.Range("M6").Value = 61.156892
$ws.Range("N6").Value = 183.470676
$ws.Range("O6").Value = 0.9308124486389074
$ws.Range("P6").Value = 0.9308124486389074
$ws.Range("Q6").Value = 58.66008034157733
$ws.Range("R6").Value = 527.940723074196
$ws.Range("S6").Value = 0.3041918857045371
$ws.Range("T6").Value = 0.3041918857045371
# Row 7
$ws.Range("G7").Value = 0.9591736666666666
$ws.Range("I7").Value = 0.3268025542087943
$ws.Range("J7").Value = 0.3268025542087943
$ws.Range("O7").Value = 0.02171808228502914
$ws.Range("P7").Value = 0.02171808228502914
$ws.Range("S7").Value = 0.007097524763264291
$ws.Range("T7").Value = 0.007097524763264291
# Row 8
$ws.Range("G8").Value = 0.9591736666666666
$ws.Range("I8").Value = 0.3268025542087943
$ws.Range("J8").Value = 0.3268025542087943
$ws.Range("M8").Value = 2.00294
$ws.Range("N8").Value = 6.00882
$ws.Range("O8").Value = 0.03048489589491914
$ws.Range("P8").Value = 0.03048489589491914
$ws.Range("Q8").Value = 1.921167303913333
$ws.Range("R8").Value = 17.29050573522
$ws.Range("S8").Value = 0.009962541843248764
$ws.Range("T8").Value = 0.009962541843248764
# Row 9
$ws.Range("G9").Value = 0.9591736666666666
$ws.Range("I9").Value = 0.3268025542087943
$ws.Range("J9").Value = 0.3268025542087943
$ws.Range("M9").Value = 1.115932333333334
$ws.Range("N9").Value = 3.347797
$ws.Range("O9").Value = 0.01698457318114416
$ws.Range("P9").Value = 0.01698457318114415
$ws.Range("Q9").Value = 1.070372907915222
$ws.Range("R9").Value = 9.633356171237
$ws.Range("S9").Value = 0.005550601897744097
$ws.Range("T9").Value = 0.005550601897744097
# Row 10
$ws.Range("I10").Value = 0.04674417878325851
$ws.Range("J10").Value = 0.04674417878325852
$ws.Range("M10").Value = 61.156892
$ws.Range("N10").Value = 183.470676
$ws.Range("O10").Value = 0.9308124486389074
$ws.Range("P10").Value = 0.9308124486389074
$ws.Range("Q10").Value = 8.390440183570666
$ws.Range("R10").Value = 75.513961652136
$ws.Range("S10").Value = 0.04351006351285972
$ws.Range("T10").Value = 0.04351006351285973
# Row 11
$ws.Range("I11").Value = 0.04674417878325851
$ws.Range("J11").Value = 0.04674417878325852
$ws.Range("O11").Value = 0.02171808228502914
$ws.Range("P11").Value = 0.02171808228502914
$ws.Range("Q11").Value = 0.1957690516288889
$ws.Range("S11").Value = 0.001015193921160922
$ws.Range("T11").Value = 0.001015193921160922
# Row 12
$ws.Range("I12").Value = 0.04674417878325851
$ws.Range("J12").Value = 0.04674417878325852
$ws.Range("M12").Value = 2.00294
$ws.Range("N12").Value = 6.00882
$ws.Range("O12").Value = 0.03048489589491914
$ws.Range("P12").Value = 0.03048489589491914
$ws.Range("Q12").Value = 0.2747940209466667
$ws.Range("R12").Value = 2.47314618852
$ws.Range("S12").Value = 0.001424991423901124
$ws.Range("T12").Value = 0.001424991423901124
# Row 13
$ws.Range("I13").Value = 0.04674417878325851
$ws.Range("J13").Value = 0.04674417878325852
$ws.Range("M13").Value = 1.115932333333334
$ws.Range("N13").Value = 3.347797
$ws.Range("O13").Value = 0.01698457318114416
$ws.Range("P13").Value = 0.01698457318114415
$ws.Range("Q13").Value = 0.1531007084491111
$ws.Range("R13").Value = 1.377906376042
$ws.Range("S13").Value = 0.0007939299253367403
$ws.Range("T13").Value = 0.0007939299253367403
# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.2294253333333333
$ws.Range("H14").Value = 0.688276
$ws.Range("I14").Value = 0.07816810191849585
$ws.Range("J14").Value = 0.07816810191849587
$ws.Range("M14").Value = 61.156892
$ws.Range("N14").Value = 183.470676
$ws.Range("O14").Value = 0.9308124486389074
$ws.Range("P14").Value = 0.9308124486389074
$ws.Range("Q14").Value = 14.03094033273067
$ws.Range("R14").Value = 126.278462994576
$ws.Range("S14").Value = 0.07275984235221081
$ws.Range("T14").Value = 0.07275984235221082
# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.2294253333333333
$ws.Range("H15").Value = 0.688276
$ws.Range("I15").Value = 0.07816810191849585
$ws.Range("J15").Value = 0.07816810191849587
$ws.Range("O15").Value = 0.02171808228502914
$ws.Range("P15").Value = 0.02171808228502914
$ws.Range("Q15").Value = 0.3273754203955556
$ws.Range("R15").Value = 2.94637878356
$ws.Range("S15").Value = 0.001697661269530437
$ws.Range("T15").Value = 0.001697661269530437
# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.2294253333333333
$ws.Range("H16").Value = 0.688276
$ws.Range("I16").Value = 0.07816810191849585
$ws.Range("J16").Value = 0.07816810191849587
$ws.Range("M16").Value = 2.00294
$ws.Range("N16").Value = 6.00882
$ws.Range("O16").Value = 0.03048489589491914
$ws.Range("P16").Value = 0.03048489589491914
$ws.Range("Q16").Value = 0.4595251771466667
$ws.Range("R16").Value = 4.13572659432
$ws.Range("S16").Value = 0.002382946449288775
$ws.Range("T16").Value = 0.002382946449288776
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.2294253333333333
$ws.Range("H17").Value = 0.688276
$ws.Range("I17").Value = 0.07816810191849585
$ws.Range("J17").Value = 0.07816810191849587
$ws.Range("M17").Value = 1.115932333333334
$ws.Range("N17").Value = 3.347797
$ws.Range("O17").Value = 0.01698457318114416
$ws.Range("P17").Value = 0.01698457318114415
$ws.Range("Q17").Value = 0.2560231475524445
$ws.Range("R17").Value = 2.304208327972
$ws.Range("S17").Value = 0.001327651847465828
$ws.Range("T17").Value = 0.001327651847465828
